$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.694.42"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "'3.741.35"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'612.65"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'178.75"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("D7").Value = "'3.739.36"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'6.58"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "'39.88"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "'4.362.94"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "'3.738.60"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "'69.762.91"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").Value = "'502.38"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'16.33"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'9.13"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +8.96%  "
$ws.Range("D25").Value = "'86.01"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").Value = "'11.77"
$ws.Range("E26").Value = "  +7.84%  "
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'0.0000136"
$ws.Range("E28").Value = "  +9.48%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").Value = "'8.08"
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("D33").Value = "'30.37"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").Value = "'0.113"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.356"
$ws.Range("E38").Value = "  +5.87%  "
$ws.Range("E39").Value = "  +3.72%  "
$ws.Range("E40").Value = "  +14.09%  "
$ws.Range("B41").Value = "Arweave"
$ws.Range("C41").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D41").Value = "'46.61"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'449.15"
$ws.Range("E42").Value = "  +8.20%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.08"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").Value = "'49.71"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").Value = "'8.55"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").Value = "'2.955.29"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'138.81"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'27.09"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("E51").Value = "  +0.62%  "
